$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update summary header values ---
$ws.Range("E11").Value = 68400      # VALOR MORA total
$ws.Range("C13").Value = 1          # Cant. Trabajadores
$ws.Range("F13").Value = 2          # Cant. Periodos

# --- Reorder the two remaining period rows (1904 before 1905) and update amounts ---
$ws.Range("E16").Value = "1904"
$ws.Range("F16").Value = 32400
$ws.Range("E17").Value = "1905"
$ws.Range("F17").Value = 36000

# --- Row 17 becomes the new last row of the table, so give it the closing-border
#     formatting that row 18 (the row about to be removed) currently has ---
$ws.Range("B18:J18").Copy()
$ws.Range("B17:J17").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Remove the now-obsolete worker row (ALDAIR TURIZO RUIZ / period 2503) ---
$ws.Rows("18:18").Delete()

$wb.Save()
